$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(17,18,19,24,25,28,29,30,31,32,33,34,35,38,39,44,49,50,65,66,71,72,73,74,75,76,91,92,99,100,103,104,109,110,119,120,127,130,131,147,148,224,225,256,257)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 5).Value = 0
}
